# Create_activity.xlsx — "forgot password pin automation added and input
# time error in create activity is solved"
#
# 1. The create-activity test data no longer drives a time_from/time_to
#    pair (they produced bad input-time errors), so those two columns are
#    removed entirely. file_path_photo / file_path_video slide left from
#    F:G into D:E.
# 2. A forgot-password-pin automation value replaces the old dummy
#    "Password" test value (12345) with a real looking secret, which Excel
#    auto-links (same as the e-mail cell right next to it).
# 3. The active selection is left on B3, matching the author's saved view.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Drop the time_from / time_to columns -----------------------------
# This removes headers D1:E1 ("time_from"/"time_to") and their row-2 time
# values, shifting file_path_photo/file_path_video from F:G into D:E.
$ws.Range("D1:E1").EntireColumn.Delete()

# --- 2. New forgot-password-pin value in the Password column -------------
$ws.Range("B2").Value = "Lakhera@1998"
$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:Lakhera@1998") | Out-Null

# --- 3. Leave the selection where the author last left it ----------------
$ws.Range("B3").Select() | Out-Null
